# Update "想去人数" (interest count) values in column F across sheets,
# reflecting the site re-scrape recorded in the gh-pages data commit.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 7761
$ws.Cells.Item(3, 6).Value = 7761
$ws.Cells.Item(5, 6).Value = 7918
$ws.Cells.Item(9, 6).Value = 6743
$ws.Cells.Item(10, 6).Value = 3402
$ws.Cells.Item(12, 6).Value = 3741
$ws.Cells.Item(13, 6).Value = 48
$ws.Cells.Item(15, 6).Value = 50
$ws.Cells.Item(16, 6).Value = 75
$ws.Cells.Item(18, 6).Value = 476
$ws.Cells.Item(20, 6).Value = 56
$ws.Cells.Item(21, 6).Value = 329
$ws.Cells.Item(23, 6).Value = 336
$ws.Cells.Item(24, 6).Value = 3890
$ws.Cells.Item(26, 6).Value = 379
$ws.Cells.Item(28, 6).Value = 293
$ws.Cells.Item(29, 6).Value = 1515
$ws.Cells.Item(32, 6).Value = 2780
$ws.Cells.Item(33, 6).Value = 1920
$ws.Cells.Item(35, 6).Value = 54
$ws.Cells.Item(37, 6).Value = 65
$ws.Cells.Item(38, 6).Value = 3760
$ws.Cells.Item(42, 6).Value = 928
$ws.Cells.Item(43, 6).Value = 558
$ws.Cells.Item(45, 6).Value = 1448
$ws.Cells.Item(46, 6).Value = 248
$ws.Cells.Item(47, 6).Value = 5
$ws.Cells.Item(49, 6).Value = 653

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(6, 6).Value = 420
$ws.Cells.Item(17, 6).Value = 180

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(5, 6).Value = 7761
$ws.Cells.Item(6, 6).Value = 7761
$ws.Cells.Item(7, 6).Value = 7918
$ws.Cells.Item(10, 6).Value = 6743
$ws.Cells.Item(11, 6).Value = 3402
$ws.Cells.Item(12, 6).Value = 3741
$ws.Cells.Item(14, 6).Value = 50
$ws.Cells.Item(15, 6).Value = 75
$ws.Cells.Item(17, 6).Value = 476
$ws.Cells.Item(19, 6).Value = 56
$ws.Cells.Item(20, 6).Value = 329
$ws.Cells.Item(22, 6).Value = 336
$ws.Cells.Item(23, 6).Value = 3890
$ws.Cells.Item(27, 6).Value = 379
$ws.Cells.Item(29, 6).Value = 293
$ws.Cells.Item(30, 6).Value = 1515
$ws.Cells.Item(33, 6).Value = 2780
$ws.Cells.Item(34, 6).Value = 1920
$ws.Cells.Item(36, 6).Value = 54
$ws.Cells.Item(39, 6).Value = 3760
$ws.Cells.Item(43, 6).Value = 928
$ws.Cells.Item(44, 6).Value = 558
$ws.Cells.Item(45, 6).Value = 180
$ws.Cells.Item(46, 6).Value = 1448
$ws.Cells.Item(47, 6).Value = 248
$ws.Cells.Item(50, 6).Value = 653
